$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - SEQUELIZE
$ws.Range("B2").Value = "17/03/2023"
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = "Existen 7 (100%) elementos de JavaScript y 0 (0%) elementos de TypeScript, lenguajes diferentes 0 (0%)"
$ws.Range("F2").Value = 7

# Row 3 - BOOKSHELF
$ws.Range("B3").Value = "17/03/2023"
$ws.Range("C3").Value = 9
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = "Existen 6 (66%) elementos de JavaScript y 0 (0%) elementos de TypeScript, lenguajes diferentes 3 (33%)"
$ws.Range("F3").Value = 9

# Row 4 - PRISMA
$ws.Range("B4").Value = "17/03/2023"
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = "Existen 2 (33%) elementos de JavaScript y 3 (50%) elementos de TypeScript, lenguajes diferentes 1 (16%)"
$ws.Range("F4").Value = 6
